$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the three "recommitment" columns right after "commitment end date" (col AB),
#        i.e. before the existing "Currency" column (col AC). Inserting three times at AC
#        shifts Currency..USD MSRP three columns to the right.
$ws.Columns("AC").Insert()
$ws.Columns("AC").Insert()
$ws.Columns("AC").Insert()
$ws.Range("AC1").Value = "recommitment"
$ws.Range("AD1").Value = "recommitment start date"
$ws.Range("AE1").Value = "recommitment end date"

# --- 2. Insert the "external reference id" column right after "Customer Name" (col Q),
#        i.e. before the existing "Seamless Move" column (col R).
$ws.Columns("R").Insert()
$ws.Range("R1").Value = "external reference id"

# --- 3. Refresh the AutoFilter so its range covers the new last column (A1:AN1).
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:AN1").AutoFilter()

# --- 4. Update the _xlnm._FilterDatabase defined name to match the new range.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Data!`$A`$1:`$AN`$1"
    }
}

# --- 5. Column widths close to the originals for the newly inserted / shifted columns
#        (the header style "s=1" is already inherited automatically by the Insert calls).
# Approximate column widths to match the new layout (best effort; Excel quantizes these).
$ws.Columns("Q").ColumnWidth = 18.17
$ws.Columns("R").ColumnWidth = 21.67
$ws.Columns("AD").ColumnWidth = 23.17
$ws.Columns("AE").ColumnWidth = 23.17
$ws.Columns("AF").ColumnWidth = 23.17

# --- 6. Update the view: scroll / selection to match the new state.
$ws.Range("R28").Select()

Write-Host "done"
